$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Edit 1: "Most (__%) ASVs received a family-level taxonomic
# assignment, so we chose to co" -> "We chose to co"
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Most (__%) ASVs received a family-level taxonomic assignment, so we chose to co",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We chose to co", 2) | Out-Null

# ---------------------------------------------------------------
# Edit 2: citation / methods text update after "Kartzinel"
# ---------------------------------------------------------------
$old2 = " " + [char]0x2013 + " find other metabarcoding ones) by summing the cumulative read abundances across the ASVs that corresponded to each diet family in each sample.  All DNA matching any predator family present on an individual sequencing run was also removed as a conservative method to account for potential sequence jumping within sequencing runs which could alter prey identity or diversity in favor of predator species on a shared run (CITE SEQUENCE JUMPING). "
$new2 = " et al., Eitzinger et al.) by summing the cumulative read abundances across the ASVs that corresponded to each diet family in each sample. Family-level data provides information comparable to previous studies, additionally, on Palmyra, each invertebrate family corresponds to an average of 1.9 (" + [char]0xB1 + " 0.13 SE) species, so for this system a family-level taxonomic assignment may closely mirror species-level assignments. All DNA matching any predator family present on an individual sequencing run was also removed as a conservative method to account for potential sequence jumping within sequencing runs which could alter prey identity or diversity in favor of predator species on a shared run (van der Valk et al. 2020)."

$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------
# Edit 3: merge "paired-end" run-split (grammar check artifact)
# back into a single run by retyping identical text
# ---------------------------------------------------------------
$old3 = " paired-end reads corresponded to "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# ---------------------------------------------------------------
# Edit 4: re-typing the "conflicting taxonomic assignments..." run
# (page-break/run split is a layout artifact that Word recalculates)
# ---------------------------------------------------------------
$old4 = " conflicting taxonomic assignments at the family level or higher between the BOLD and BLAST assignments"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $old4, 2) | Out-Null

Write-Output "done"
